# Auto-generated cell value updates based on diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.870921420471129
$ws.Range("C2").Value = 4.440271573361471
$ws.Range("E2").Value = 12.32198667410082
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 31.9548318330852
$ws.Range("H2").Value = 15.07117146919645
$ws.Range("I2").Value = 21.67611151978054
$ws.Range("K2").Value = 9.001189111014824
$ws.Range("M2").Value = 13.87158067594262
$ws.Range("N2").Value = 18.96057001254424
$ws.Range("B3").Value = 8.606907451801838
$ws.Range("C3").Value = 4.190801336075801
$ws.Range("E3").Value = 12.1006781926643
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 31.98461888728828
$ws.Range("H3").Value = 15.1190594202505
$ws.Range("I3").Value = 21.7582661203378
$ws.Range("K3").Value = 8.831632129853222
$ws.Range("M3").Value = 13.71012614330197
$ws.Range("N3").Value = 19.02664034297902
$ws.Range("B4").Value = 8.44253328065043
$ws.Range("C4").Value = 4.028974843864933
$ws.Range("E4").Value = 11.96689197690076
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 32.01369230313371
$ws.Range("H4").Value = 15.15109970769102
$ws.Range("I4").Value = 21.81306741508694
$ws.Range("K4").Value = 8.727837983983655
$ws.Range("M4").Value = 13.61362722133091
$ws.Range("N4").Value = 19.06905943395833
$ws.Range("B5").Value = 8.375089973111404
$ws.Range("C5").Value = 3.960874587790245
$ws.Range("E5").Value = 11.91298060753073
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 32.02824113243551
$ws.Range("H5").Value = 15.16481862755806
$ws.Range("I5").Value = 21.83649308554101
$ws.Range("K5").Value = 8.685679897076406
$ws.Range("M5").Value = 13.57501084736904
$ws.Range("N5").Value = 19.08681259651002
$ws.Range("B6").Value = 8.363866690458353
$ws.Range("C6").Value = 3.949437166472746
$ws.Range("E6").Value = 11.90406768888711
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 32.03081975458925
$ws.Range("H6").Value = 15.1671366172839
$ws.Range("I6").Value = 21.84044889358224
$ws.Range("K6").Value = 8.678689728165418
$ws.Range("M6").Value = 13.5686426633367
$ws.Range("N6").Value = 19.08978874560153
$ws.Range("B7").Value = 8.441625425341922
$ws.Range("C7").Value = 4.028065112891142
$ws.Range("E7").Value = 11.9661623439039
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 32.01387759214375
$ws.Range("H7").Value = 15.15128204556481
$ws.Range("I7").Value = 21.81337891642343
$ws.Range("K7").Value = 8.727268783455171
$ws.Range("M7").Value = 13.61310350330114
$ws.Range("N7").Value = 19.06929696641734
$ws.Range("B8").Value = 8.780424576410187
$ws.Range("C8").Value = 4.356060506969049
$ws.Range("E8").Value = 12.24529219237578
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 31.96285837902327
$ws.Range("H8").Value = 15.08713545526357
$ws.Range("I8").Value = 21.70353261597421
$ws.Range("K8").Value = 8.942696173544499
$ws.Range("M8").Value = 13.81539183097366
$ws.Range("N8").Value = 18.98296761868388
$ws.Range("B9").Value = 9.422307578289466
$ws.Range("C9").Value = 4.929956625939745
$ws.Range("E9").Value = 12.80578991079757
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 31.94876431683474
$ws.Range("H9").Value = 14.98230345734545
$ws.Range("I9").Value = 21.52280012184904
$ws.Range("K9").Value = 9.36507163857307
$ws.Range("M9").Value = 14.23093701644065
$ws.Range("N9").Value = 18.82830151417549
$ws.Range("B10").Value = 9.874700020398587
$ws.Range("C10").Value = 5.308740668031085
$ws.Range("E10").Value = 13.22092472769101
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 31.9912141614252
$ws.Range("H10").Value = 14.9181106682417
$ws.Range("I10").Value = 21.41128731895587
$ws.Range("K10").Value = 9.672008901580229
$ws.Range("M10").Value = 14.54493089021912
$ws.Range("N10").Value = 18.72348957618331
$ws.Range("B11").Value = 10.07533524797531
$ws.Range("C11").Value = 5.471685562748968
$ws.Range("E11").Value = 13.40955049165066
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 32.02203535977183
$ws.Range("H11").Value = 14.89170359383026
$ws.Range("I11").Value = 21.3652044661152
$ws.Range("K11").Value = 9.810243249125389
$ws.Range("M11").Value = 14.68905907794106
$ws.Range("N11").Value = 18.67770307797397
$ws.Range("B12").Value = 10.15049564455571
$ws.Range("C12").Value = 5.532038688596482
$ws.Range("E12").Value = 13.48087036255887
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 32.03536174458659
$ws.Range("H12").Value = 14.88210665395329
$ws.Range("I12").Value = 21.34842440026158
$ws.Range("K12").Value = 9.862336682606488
$ws.Range("M12").Value = 14.74377202557997
$ws.Range("N12").Value = 18.66063563181177
$ws.Range("B13").Value = 10.13434588695462
$ws.Range("C13").Value = 5.519100671334506
$ws.Range("E13").Value = 13.46551637049609
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 32.0324180853912
$ws.Range("H13").Value = 14.88415559433387
$ws.Range("I13").Value = 21.35200842557335
$ws.Range("K13").Value = 9.851129485384364
$ws.Range("M13").Value = 14.73198342429977
$ws.Range("N13").Value = 18.66429938365529
$ws.Range("B14").Value = 10.08153548504389
$ws.Range("C14").Value = 5.476677933460897
$ws.Range("E14").Value = 13.41542062529051
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 32.02309856813785
$ws.Range("H14").Value = 14.89090596854501
$ws.Range("I14").Value = 21.363810507232
$ws.Range("K14").Value = 9.814534367239677
$ws.Range("M14").Value = 14.69355793399462
$ws.Range("N14").Value = 18.67629350694149
$ws.Range("B15").Value = 10.04907928501356
$ws.Range("C15").Value = 5.450516816126616
$ws.Range("E15").Value = 13.38471916697882
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 32.01760559133028
$ws.Range("H15").Value = 14.89509325931277
$ws.Range("I15").Value = 21.37112701966856
$ws.Range("K15").Value = 9.792084350979934
$ws.Range("M15").Value = 14.67003725838905
$ws.Range("N15").Value = 18.68367549543469
$ws.Range("B16").Value = 9.861477720369841
$ws.Range("C16").Value = 5.297903006911531
$ws.Range("E16").Value = 13.20858650026799
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 31.98943158406126
$ws.Range("H16").Value = 14.91989275140764
$ws.Range("I16").Value = 21.41439262967437
$ws.Range("K16").Value = 9.66294237060697
$ws.Range("M16").Value = 14.53553338497178
$ws.Range("N16").Value = 18.72651981162937
$ws.Range("B17").Value = 9.745013671030962
$ws.Range("C17").Value = 5.201877282372823
$ws.Range("E17").Value = 13.10042250728402
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 31.9750964747441
$ws.Range("H17").Value = 14.93582289935346
$ws.Range("I17").Value = 21.44212618678992
$ws.Range("K17").Value = 9.583322860399697
$ws.Range("M17").Value = 14.45331383726563
$ws.Range("N17").Value = 18.75328732001956
$ws.Range("B18").Value = 9.67754483305111
$ws.Range("C18").Value = 5.145764361786235
$ws.Range("E18").Value = 13.03819302784176
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 31.96793487445734
$ws.Range("H18").Value = 14.94524847664447
$ws.Range("I18").Value = 21.45851483967664
$ws.Range("K18").Value = 9.537399686410884
$ws.Range("M18").Value = 14.40614784180427
$ws.Range("N18").Value = 18.76886153596048
$ws.Range("B19").Value = 9.654620665382406
$ws.Range("C19").Value = 5.126614243829431
$ws.Range("E19").Value = 13.01712276983688
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 31.96569615285523
$ws.Range("H19").Value = 14.94848495128789
$ws.Range("I19").Value = 21.46413873988303
$ws.Range("K19").Value = 9.521830530022916
$ws.Range("M19").Value = 14.39020119146116
$ws.Range("N19").Value = 18.77416535285337
$ws.Range("B20").Value = 9.757461890148852
$ws.Range("C20").Value = 5.212190645858293
$ws.Range("E20").Value = 13.11193901429374
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 31.97651031953264
$ws.Range("H20").Value = 14.93409988345104
$ws.Range("I20").Value = 21.43912865106198
$ws.Range("K20").Value = 9.591812144985431
$ws.Range("M20").Value = 14.46205373230363
$ws.Range("N20").Value = 18.75041943396646
$ws.Range("B21").Value = 10.09706987190336
$ws.Range("C21").Value = 5.489175195246342
$ws.Range("E21").Value = 13.43013850342129
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 32.02579103030895
$ws.Range("H21").Value = 14.88891227904474
$ws.Range("I21").Value = 21.36032573179125
$ws.Range("K21").Value = 9.825290509679277
$ws.Range("M21").Value = 14.70484117766107
$ws.Range("N21").Value = 18.67276320016653
$ws.Range("B22").Value = 10.31423912839042
$ws.Range("C22").Value = 5.662327893341002
$ws.Range("E22").Value = 13.63743397375628
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 32.06764418673564
$ws.Range("H22").Value = 14.86172811140488
$ws.Range("I22").Value = 21.31273290389838
$ws.Range("K22").Value = 9.976385576747543
$ws.Range("M22").Value = 14.86427904528561
$ws.Range("N22").Value = 18.62358881941919
$ws.Range("B23").Value = 10.19879201997297
$ws.Range("C23").Value = 5.570634195711619
$ws.Range("E23").Value = 13.52688194880871
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 32.04442443563521
$ws.Range("H23").Value = 14.87602159755586
$ws.Range("I23").Value = 21.33777552890328
$ws.Range("K23").Value = 9.895896867988661
$ws.Range("M23").Value = 14.7791309373424
$ws.Range("N23").Value = 18.64969011007286
$ws.Range("B24").Value = 9.751835639683891
$ws.Range("C24").Value = 5.207530797376331
$ws.Range("E24").Value = 13.10673253440134
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 31.97586775660572
$ws.Range("H24").Value = 14.93487802691051
$ws.Range("I24").Value = 21.44048245326716
$ws.Range("K24").Value = 9.587974599029888
$ws.Range("M24").Value = 14.45810210066228
$ws.Range("N24").Value = 18.75171542817158
$ws.Range("B25").Value = 9.251666839197819
$ws.Range("C25").Value = 4.782197271220864
$ws.Range("E25").Value = 12.65324251033029
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 31.94332526726306
$ws.Range("H25").Value = 15.00841425476486
$ws.Range("I25").Value = 21.56796711504038
$ws.Range("K25").Value = 9.251168032886321
$ws.Range("M25").Value = 14.11678954883302
$ws.Range("N25").Value = 18.86858687676363
